# The "Input" column (D) on the Channels sheet listed snake/bus input jacks
# labelled "B1".."B16". Relabel them to "A1".."A16" (16 cells, rows 17-35,
# skipping the blank "Unused" rows 19/25/34).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Channels")

$rows = @(17, 18, 20, 21, 22, 23, 24, 26, 27, 28, 29, 30, 31, 32, 33, 35)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $newValue = "A" + ($i + 1)
    $ws.Range("D" + $rows[$i]).Value = $newValue
}
